$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues = -4163
$xlPasteValues = -4163

# --- Step 1: write every changed cell as a formula that evaluates to the
#     literal text (so numeric-looking strings like "0.999" or "5.34" stay
#     text, matching the original inlineStr cells) ---
$ws.Range("D2").Formula = '="58.275.99"'
$ws.Range("E2").Formula = '="  -1.19%  "'
$ws.Range("D3").Formula = '="2.475.97"'
$ws.Range("E3").Formula = '="  -1.06%  "'
$ws.Range("D4").Formula = '="0.999"'
$ws.Range("E4").Formula = '="  -0.10%  "'
$ws.Range("D5").Formula = '="521.23"'
$ws.Range("E5").Formula = '="  -2.33%  "'
$ws.Range("D6").Formula = '="134.58"'
$ws.Range("E6").Formula = '="  -1.02%  "'
$ws.Range("D7").Formula = '="0.997"'
$ws.Range("E7").Formula = '="  -0.22%  "'
$ws.Range("E8").Formula = '="  -1.51%  "'
$ws.Range("D9").Formula = '="2.488.38"'
$ws.Range("E9").Formula = '="  -0.65%  "'
$ws.Range("D10").Formula = '="0.0985"'
$ws.Range("E10").Formula = '="  -3.38%  "'
$ws.Range("E11").Formula = '="  -1.02%  "'
$ws.Range("D12").Formula = '="5.34"'
$ws.Range("E12").Formula = '="  -1.04%  "'
$ws.Range("D13").Formula = '="0.339"'
$ws.Range("E13").Formula = '="  -2.37%  "'
$ws.Range("D14").Formula = '="2.914.16"'
$ws.Range("E14").Formula = '="  -1.08%  "'
$ws.Range("D15").Formula = '="58.208.23"'
$ws.Range("E15").Formula = '="  -1.19%  "'
$ws.Range("D16").Formula = '="22.19"'
$ws.Range("E16").Formula = '="  -2.30%  "'
$ws.Range("E17").Formula = '="  -2.13%  "'
$ws.Range("D18").Formula = '="2.473.39"'
$ws.Range("E18").Formula = '="  -0.73%  "'
$ws.Range("D19").Formula = '="10.69"'
$ws.Range("E19").Formula = '="  -3.07%  "'
$ws.Range("D20").Formula = '="4.19"'
$ws.Range("E20").Formula = '="  -1.61%  "'
$ws.Range("D21").Formula = '="320.41"'
$ws.Range("E21").Formula = '="  -1.17%  "'
$ws.Range("E22").Formula = '="  +0.03%  "'
$ws.Range("E23").Formula = '="  -5.45%  "'
$ws.Range("D24").Formula = '="64.53"'
$ws.Range("E24").Formula = '="  -1.03%  "'
$ws.Range("D25").Formula = '="0.411"'
$ws.Range("E25").Formula = '="  -2.64%  "'
$ws.Range("E26").Formula = '="  -0.29%  "'
$ws.Range("E27").Formula = '="  -1.94%  "'
$ws.Range("D28").Formula = '="7.40"'
$ws.Range("E28").Formula = '="  -1.76%  "'
$ws.Range("D29").Formula = '="0.0₃0751"'
$ws.Range("E29").Formula = '="  -1.32%  "'
$ws.Range("D30").Formula = '="169.19"'
$ws.Range("E30").Formula = '="  -0.63%  "'
$ws.Range("E31").Formula = '="  -2.73%  "'
$ws.Range("B32").Formula = '="Fetch.AI"'
$ws.Range("C32").Formula = '="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"'
$ws.Range("D32").Formula = '="1.19"'
$ws.Range("E32").Formula = '="  +1.67%  "'
$ws.Range("B33").Formula = '="Aptos"'
$ws.Range("C33").Formula = '="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"'
$ws.Range("D33").Formula = '="6.31"'
$ws.Range("E33").Formula = '="  -2.14%  "'
$ws.Range("E34").Formula = '="  -0.01%  "'
$ws.Range("E35").Formula = '="  -0.33%  "'
$ws.Range("E36").Formula = '="  -1.57%  "'
$ws.Range("D37").Formula = '="1.33"'
$ws.Range("E37").Formula = '="  -1.74%  "'
$ws.Range("D38").Formula = '="4.01"'
$ws.Range("E38").Formula = '="  -0.99%  "'
$ws.Range("D39").Formula = '="36.60"'
$ws.Range("E39").Formula = '="  -0.39%  "'
$ws.Range("E40").Formula = '="  -3.18%  "'
$ws.Range("E41").Formula = '="  +0.20%  "'
$ws.Range("D42").Formula = '="5.17"'
$ws.Range("E42").Formula = '="  +3.08%  "'
$ws.Range("D43").Formula = '="3.46"'
$ws.Range("E43").Formula = '="  -3.24%  "'
$ws.Range("D44").Formula = '="274.95"'
$ws.Range("E44").Formula = '="  -2.44%  "'
$ws.Range("D45").Formula = '="0.598"'
$ws.Range("E45").Formula = '="  -0.36%  "'
$ws.Range("D46").Formula = '="124.25"'
$ws.Range("E46").Formula = '="  -4.41%  "'
$ws.Range("D47").Formula = '="0.0911"'
$ws.Range("E47").Formula = '="  -1.58%  "'
$ws.Range("D48").Formula = '="0.0491"'
$ws.Range("E48").Formula = '="  -1.69%  "'
$ws.Range("D49").Formula = '="0.0214"'
$ws.Range("E49").Formula = '="  -2.23%  "'
$ws.Range("D50").Formula = '="17.05"'
$ws.Range("E50").Formula = '="  -1.10%  "'
$ws.Range("D51").Formula = '="1.741.08"'
$ws.Range("E51").Formula = '="  -0.96%  "'

# --- Step 2: convert those formulas back into plain static text values,
#     via copy / paste-special (values only), in contiguous blocks ---
$rng = $ws.Range("B32:C33")
$rng.Copy()
$rng.PasteSpecial($xlPasteValues)
$rng = $ws.Range("D2:D7")
$rng.Copy()
$rng.PasteSpecial($xlPasteValues)
$rng = $ws.Range("D9:D10")
$rng.Copy()
$rng.PasteSpecial($xlPasteValues)
$rng = $ws.Range("D12:D16")
$rng.Copy()
$rng.PasteSpecial($xlPasteValues)
$rng = $ws.Range("D18:D21")
$rng.Copy()
$rng.PasteSpecial($xlPasteValues)
$rng = $ws.Range("D24:D25")
$rng.Copy()
$rng.PasteSpecial($xlPasteValues)
$rng = $ws.Range("D28:D30")
$rng.Copy()
$rng.PasteSpecial($xlPasteValues)
$rng = $ws.Range("D32:D33")
$rng.Copy()
$rng.PasteSpecial($xlPasteValues)
$rng = $ws.Range("D37:D39")
$rng.Copy()
$rng.PasteSpecial($xlPasteValues)
$rng = $ws.Range("D42:D51")
$rng.Copy()
$rng.PasteSpecial($xlPasteValues)
$rng = $ws.Range("E2:E51")
$rng.Copy()
$rng.PasteSpecial($xlPasteValues)

$excel.CutCopyMode = $false
Write-Output "Applied cryptos update"
